$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.021.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.260.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.06"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.87"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.02%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.822.57"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.04%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.43"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.034.13"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.261.08"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.34"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.45"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.40"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.83"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.534"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.08"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.63"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0913"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.93"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.71"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.25"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.43"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.51%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.67"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.35"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.49"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.840.26"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.02%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0318"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +8.75%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0718"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.73"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.31"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.11"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.732"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.301.70"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.94%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.69%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.30"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.26"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.805"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +7.67%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.09%  "
